$wb = $excel.ActiveWorkbook

# --- Sheet 1: Metadata ---
$ws1 = $wb.Worksheets.Item(1)

# Version: 5.0.0 -> 6.0.0
$ws1.Cells.Item(3, 2).Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws1.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty -> "Alvearie Team"
$ws1.Cells.Item(9, 2).Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> "Jurisdiction" / "United States of America"
$ws1.Cells.Item(10, 1).Value = "Jurisdiction"
$ws1.Cells.Item(10, 2).Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row -> remove it entirely
$ws1.Rows.Item(11).Delete()

# --- Sheet 2: Elements ---
$ws2 = $wb.Worksheets.Item(2)

# Row 2 (the Extension root element) gets a custom Short/Definition instead of the generic text
$ws2.Cells.Item(2, 11).Value = "Claim Group Identifier"
$ws2.Cells.Item(2, 12).Value = "Group ID of the plan member. This ID associates the primary plan holder with any dependents."
